# Update cryptos list values per Tue Nov 14 18:57:58 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.269.91'
$ws.Range("E2").Value = '  -4.03%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.961.68'
$ws.Range("E3").Value = '  -6.16%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.73%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.70'
$ws.Range("E5").Value = '  -3.89%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.619'
$ws.Range("E6").Value = '  -5.23%  '

$ws.Range("E7").Value = '  +0.45%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '53.39'
$ws.Range("E8").Value = '  -1.38%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '61.10'
$ws.Range("E9").Value = '  +3.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.346'
$ws.Range("E10").Value = '  -5.30%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0706'
$ws.Range("E11").Value = '  -7.13%  '

$ws.Range("E12").Value = '  -5.31%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.858'
$ws.Range("E13").Value = '  -5.09%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.305.94'
$ws.Range("E14").Value = '  -3.80%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.62'
$ws.Range("E15").Value = '  -9.29%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.08'
$ws.Range("E16").Value = '  -7.68%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.996.37'
$ws.Range("E17").Value = '  -4.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '35.417.71'
$ws.Range("E18").Value = '  -3.57%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.57'
$ws.Range("E19").Value = '  -3.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.93'
$ws.Range("E20").Value = '  -5.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0813'
$ws.Range("E21").Value = '  -7.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.82'
$ws.Range("E22").Value = '  -4.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.90'
$ws.Range("E23").Value = '  -9.95%  '

$ws.Range("E24").Value = '  -0.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.29'
$ws.Range("E25").Value = '  -4.41%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.18'
$ws.Range("E26").Value = '  +0.86%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.00'
$ws.Range("E27").Value = '  -4.19%  '

$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.80'
$ws.Range("E28").Value = '  -9.47%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.09'
$ws.Range("E29").Value = '  -7.49%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.117'
$ws.Range("E30").Value = '  -4.92%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.11'
$ws.Range("E31").Value = '  -4.54%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.66'
$ws.Range("E32").Value = '  -12.09%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0573'
$ws.Range("E33").Value = '  -5.99%  '

$ws.Range("B34").Value = 'BinanceUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.01'
$ws.Range("E34").Value = '  +0.64%  '

$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.16'
$ws.Range("E35").Value = '  -11.86%  '

$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0866'
$ws.Range("E36").Value = '  +4.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.81'
$ws.Range("E37").Value = '  -1.64%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.08'
$ws.Range("E38").Value = '  -14.16%  '

$ws.Range("B39").Value = 'HuobiToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.89'
$ws.Range("E39").Value = '  +0.81%  '

$ws.Range("B40").Value = 'THORChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.70'
$ws.Range("E40").Value = '  -3.43%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.14'
$ws.Range("E41").Value = '  -10.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0204'
$ws.Range("E42").Value = '  -7.36%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.04'
$ws.Range("E43").Value = '  -9.89%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.349.55'
$ws.Range("E44").Value = '  -2.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0865'
$ws.Range("E45").Value = '  -9.22%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '87.63'
$ws.Range("E46").Value = '  -8.87%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.13'
$ws.Range("E47").Value = '  -3.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '14.90'
$ws.Range("E48").Value = '  -7.24%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.85'
$ws.Range("E49").Value = '  -1.89%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.187.97'
$ws.Range("E50").Value = '  -4.21%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.70'
$ws.Range("E51").Value = '  -4.17%  '

